# Updated cryptos list -- refreshed Price / Volume(1h) figures, and
# restored the correct Coin/Link pairing for rows 45-46 (EnergySwap <-> Decentraland).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (e.g. "1.003" or "331.80") -- force Text format first so the literal string is kept,
# then restore the default "Normal" style so no stray formatting is left behind.
$textCells = @(
    'D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D14', 'D15', 'D17', 'D18',
    'D19', 'D20', 'D21', 'D23', 'D24', 'D25', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33',
    'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46',
    'D47', 'D48', 'D49', 'D50', 'D51'
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = '@'
}

# --- Price column updates (column D) ---
$ws.Range('D2').Value = '27.483.76'
$ws.Range('D3').Value = '1.830.49'
$ws.Range('D4').Value = '1.003'
$ws.Range('D5').Value = '331.80'
$ws.Range('D6').Value = '1.003'
$ws.Range('D7').Value = '0.4581'
$ws.Range('D8').Value = '0.3806'
$ws.Range('D9').Value = '46.39'
$ws.Range('D10').Value = '0.07911'
$ws.Range('D11').Value = '0.9703'
$ws.Range('D12').Value = '21.07'
$ws.Range('D13').Value = '1.832.26'
$ws.Range('D14').Value = '5.881'
$ws.Range('D15').Value = '7.071'
$ws.Range('D17').Value = '89.53'
$ws.Range('D18').Value = '0.06638'
$ws.Range('D19').Value = '0.00001026'
$ws.Range('D20').Value = '17.12'
$ws.Range('D21').Value = '1.003'
$ws.Range('D22').Value = '27.453.40'
$ws.Range('D23').Value = '5.337'
$ws.Range('D24').Value = '10.82'
$ws.Range('D25').Value = '2.292'
$ws.Range('D26').Value = '2.031.06'
$ws.Range('D27').Value = '155.64'
$ws.Range('D28').Value = '19.40'
$ws.Range('D29').Value = '2.066'
$ws.Range('D30').Value = '5.298'
$ws.Range('D31').Value = '118.50'
$ws.Range('D32').Value = '0.9423'
$ws.Range('D33').Value = '0.09302'
$ws.Range('D34').Value = '3.585'
$ws.Range('D35').Value = '5.251'
$ws.Range('D36').Value = '1.331'
$ws.Range('D37').Value = '0.05941'
$ws.Range('D38').Value = '0.02179'
$ws.Range('D39').Value = '8.076'
$ws.Range('D40').Value = '1.144'
$ws.Range('D41').Value = '0.5774'
$ws.Range('D42').Value = '0.1828'
$ws.Range('D43').Value = '9.984'
$ws.Range('D44').Value = '1.265'
$ws.Range('D45').Value = '0.5451'
$ws.Range('D46').Value = '11.94'
$ws.Range('D47').Value = '1.871'
$ws.Range('D48').Value = '111.07'
$ws.Range('D49').Value = '0.06611'
$ws.Range('D50').Value = '1.031'
$ws.Range('D51').Value = '1.043'

# --- Volume(1h) column updates (column E) ---
$ws.Range('E2').Value = '  -0.64%  '
$ws.Range('E3').Value = '  -1.89%  '
$ws.Range('E4').Value = '  -0.70%  '
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('E6').Value = '  -0.53%  '
$ws.Range('E8').Value = '  -3.08%  '
$ws.Range('E9').Value = '  +2.13%  '
$ws.Range('E10').Value = '  -1.08%  '
$ws.Range('E11').Value = '  -3.24%  '
$ws.Range('E12').Value = '  -3.55%  '
$ws.Range('E13').Value = '  -1.94%  '
$ws.Range('E14').Value = '  -2.01%  '
$ws.Range('E15').Value = '  -2.53%  '
$ws.Range('E16').Value = '  -0.59%  '
$ws.Range('E17').Value = '  +1.15%  '
$ws.Range('E18').Value = '  -1.28%  '
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('E21').Value = '  -0.68%  '
$ws.Range('E22').Value = '  -0.67%  '
$ws.Range('E23').Value = '  -2.32%  '
$ws.Range('E24').Value = '  -0.93%  '
$ws.Range('E25').Value = '  -1.00%  '
$ws.Range('E26').Value = '  -2.64%  '
$ws.Range('E28').Value = '  -2.05%  '
$ws.Range('E29').Value = '  -4.22%  '
$ws.Range('E30').Value = '  -2.85%  '
$ws.Range('E31').Value = '  -2.76%  '
$ws.Range('E32').Value = '  -4.15%  '
$ws.Range('E33').Value = '  -2.05%  '
$ws.Range('E34').Value = '  -0.83%  '
$ws.Range('E35').Value = '  -1.24%  '
$ws.Range('E36').Value = '  -0.55%  '
$ws.Range('E37').Value = '  -2.00%  '
$ws.Range('E38').Value = '  -2.43%  '
$ws.Range('E39').Value = '  -3.25%  '
$ws.Range('E40').Value = '  -4.22%  '
$ws.Range('E41').Value = '  -3.50%  '
$ws.Range('E42').Value = '  -3.14%  '
$ws.Range('E43').Value = '  -3.05%  '
$ws.Range('E44').Value = '  +1.42%  '
$ws.Range('E45').Value = '  -3.65%  '
$ws.Range('E46').Value = '  -2.20%  '
$ws.Range('E47').Value = '  -2.78%  '
$ws.Range('E48').Value = '  -0.85%  '
$ws.Range('E49').Value = '  -2.20%  '
$ws.Range('E50').Value = '  +2.08%  '
$ws.Range('E51').Value = '  -1.40%  '

# --- Row 45 / 46: EnergySwap and Decentraland swapped position (Coin name + Link) ---
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'

# Restore default styling on the cells we protected above.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = 'Normal'
}
